# Auto-generated Excel COM-interop script applying numeric corrections
# to financial projection values (H, I, J, K, L, M, N columns) across
# several rows in multiple worksheets, per scheduled-runner sheet update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 859.1111
$ws.Range("J6").Value = 956.5
$ws.Range("L6").Value = 2869.5
$ws.Range("N6").Value = -3093.5

$ws.Range("H15").Value = 184.07
$ws.Range("I15").Value = 184.07
$ws.Range("K15").Value = 552.21
$ws.Range("M15").Value = -383.21

$ws.Range("H19").Value = 5972.9
$ws.Range("I19").Value = 9209.799999999999
$ws.Range("J19").Value = 2736
$ws.Range("K19").Value = 9209.799999999999
$ws.Range("L19").Value = 2736
$ws.Range("M19").Value = -9034.799999999999
$ws.Range("N19").Value = -3086

$ws.Range("H21").Value = 23254.25
$ws.Range("I21").Value = 17672.334
$ws.Range("K21").Value = 17672.334
$ws.Range("M21").Value = -17204.334

$ws.Range("H23").Value = 23254.25
$ws.Range("I23").Value = 17672.334
$ws.Range("K23").Value = 17672.334
$ws.Range("M23").Value = -17438.334

$ws.Range("H33").Value = 9804280
$ws.Range("I33").Value = 11494604
$ws.Range("K33").Value = 11494604
$ws.Range("M33").Value = -11494375

$ws.Range("H38").Value = 486.5
$ws.Range("I38").Value = 178.4
$ws.Range("K38").Value = 535.2
$ws.Range("M38").Value = -163.2

$ws.Range("H40").Value = 2076
$ws.Range("I40").Value = 4726.3335
$ws.Range("J40").Value = 1545.9333
$ws.Range("K40").Value = 4726.3335
$ws.Range("L40").Value = 1545.9333
$ws.Range("M40").Value = -4551.3335
$ws.Range("N40").Value = -1895.9333

$ws.Range("H64").Value = 2723.5833
$ws.Range("I64").Value = 3175
$ws.Range("J64").Value = 2497.875
$ws.Range("K64").Value = 3175
$ws.Range("L64").Value = 2497.875
$ws.Range("M64").Value = -2927
$ws.Range("N64").Value = -2993.875

$ws.Range("H67").Value = 2723.5833
$ws.Range("I67").Value = 3175
$ws.Range("J67").Value = 2497.875
$ws.Range("K67").Value = 3175
$ws.Range("L67").Value = 2497.875
$ws.Range("M67").Value = -2317
$ws.Range("N67").Value = -4213.875

$ws.Range("H116").Value = 2614.389
$ws.Range("I116").Value = 1515
$ws.Range("J116").Value = 2928.5
$ws.Range("K116").Value = 1515
$ws.Range("L116").Value = 2928.5
$ws.Range("M116").Value = 1927
$ws.Range("N116").Value = -9812.5

$ws.Range("H138").Value = 2873.3132
$ws.Range("I138").Value = 1338.4348
$ws.Range("J138").Value = 3337.8157
$ws.Range("K138").Value = 4015.3044
$ws.Range("L138").Value = 10013.4471
$ws.Range("M138").Value = 1124.6956
$ws.Range("N138").Value = -20293.4471

$ws.Range("H139").Value = 35500
$ws.Range("J139").Value = 35500
$ws.Range("L139").Value = 35500
$ws.Range("N139").Value = -45780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 9537.111000000001
$ws.Range("J37").Value = 18850
$ws.Range("L37").Value = 18850
$ws.Range("N37").Value = -19396

$ws.Range("H45").Value = 999.0769
$ws.Range("I45").Value = 887.6667
$ws.Range("J45").Value = 1249.75
$ws.Range("K45").Value = 887.6667
$ws.Range("L45").Value = 1249.75
$ws.Range("M45").Value = -510.6667
$ws.Range("N45").Value = -2003.75

$ws.Range("H132").Value = 4514.7026
$ws.Range("I132").Value = 4924.6143
$ws.Range("J132").Value = 3140.2942
$ws.Range("K132").Value = 14773.8429
$ws.Range("L132").Value = 9420.882599999999
$ws.Range("M132").Value = -12243.8429
$ws.Range("N132").Value = -14480.8826

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 39950
$ws.Range("J55").Value = 39950
$ws.Range("L55").Value = 39950
$ws.Range("N55").Value = -40496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 980
$ws.Range("K16").Value = 980
$ws.Range("M16").Value = -693

$ws.Range("H99").Value = 2073.4546
$ws.Range("I99").Value = 1722
$ws.Range("J99").Value = 2446.875
$ws.Range("K99").Value = 1722
$ws.Range("L99").Value = 2446.875
$ws.Range("M99").Value = -224
$ws.Range("N99").Value = -5442.875

$ws.Range("H102").Value = 48333.332
$ws.Range("J102").Value = 48333.332
$ws.Range("L102").Value = 48333.332
$ws.Range("N102").Value = -53201.332

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 980
$ws.Range("K113").Value = 980
$ws.Range("M113").Value = 1190

$ws.Range("H126").Value = 2073.4546
$ws.Range("I126").Value = 1722
$ws.Range("J126").Value = 2446.875
$ws.Range("K126").Value = 5166
$ws.Range("L126").Value = 7340.625
$ws.Range("M126").Value = -2696
$ws.Range("N126").Value = -12280.625

$ws.Range("H141").Value = 72342.86
$ws.Range("J141").Value = 72342.86
$ws.Range("L141").Value = 72342.86
$ws.Range("N141").Value = -82702.86

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 670.617
$ws.Range("I113").Value = 590.1053000000001
$ws.Range("J113").Value = 1010.55554
$ws.Range("K113").Value = 1770.3159
$ws.Range("L113").Value = 3031.66662
$ws.Range("M113").Value = 399.6840999999999
$ws.Range("N113").Value = -7371.66662

$ws.Range("H132").Value = 2315.0454
$ws.Range("J132").Value = 2499.9473
$ws.Range("L132").Value = 22499.5257
$ws.Range("N132").Value = -27559.5257

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1223.9025
$ws.Range("I97").Value = 1136.129
$ws.Range("J97").Value = 1496
$ws.Range("K97").Value = 1136.129
$ws.Range("L97").Value = 1496
$ws.Range("M97").Value = -640.1289999999999
$ws.Range("N97").Value = -2488

$ws.Range("H102").Value = 1713.52
$ws.Range("I102").Value = 1711.2
$ws.Range("K102").Value = 1711.2
$ws.Range("M102").Value = -89.20000000000005

$ws.Range("H122").Value = 21278902
$ws.Range("I122").Value = 32260226
$ws.Range("K122").Value = 96780678
$ws.Range("M122").Value = -96778228

$ws.Range("H132").Value = 7892.6816
$ws.Range("I132").Value = 9476.134
$ws.Range("J132").Value = 4499.5713
$ws.Range("K132").Value = 28428.402
$ws.Range("L132").Value = 13498.7139
$ws.Range("M132").Value = -25898.402
$ws.Range("N132").Value = -18558.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2944.389
$ws.Range("I40").Value = 2633.2666
$ws.Range("K40").Value = 2633.2666
$ws.Range("M40").Value = -2497.2666

$ws.Range("H55").Value = 331.82144
$ws.Range("I55").Value = 299.6111
$ws.Range("J55").Value = 389.8
$ws.Range("K55").Value = 299.6111
$ws.Range("L55").Value = 389.8
$ws.Range("M55").Value = -126.6111
$ws.Range("N55").Value = -735.8

$ws.Range("H68").Value = 1200.3334
$ws.Range("I68").Value = 1200.3334
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1200.3334
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -451.3334
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1200.3334
$ws.Range("I71").Value = 1200.3334
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6001.666999999999
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2257.666999999999
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 350.63635
$ws.Range("I100").Value = 299.75
$ws.Range("J100").Value = 486.33334
$ws.Range("K100").Value = 599.5
$ws.Range("L100").Value = 972.66668
$ws.Range("M100").Value = -58.5
$ws.Range("N100").Value = -2054.66668
